$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4400
$ws.Range("J64").Value = 4100
$ws.Range("L64").Value = 4100
$ws.Range("N64").Value = -4596

$ws.Range("H67").Value = 4400
$ws.Range("J67").Value = 4100
$ws.Range("L67").Value = 4100
$ws.Range("N67").Value = -5816

$ws.Range("H74").Value = 3911.111
$ws.Range("I74").Value = 3837.5
$ws.Range("J74").Value = 4500
$ws.Range("K74").Value = 3837.5
$ws.Range("L74").Value = 4500
$ws.Range("M74").Value = -2901.5
$ws.Range("N74").Value = -6372

$ws.Range("H76").Value = 3068.7317
$ws.Range("I76").Value = 3097.6572
$ws.Range("J76").Value = 2900
$ws.Range("K76").Value = 3097.6572
$ws.Range("L76").Value = 2900
$ws.Range("M76").Value = -2782.6572
$ws.Range("N76").Value = -3530

$ws.Range("H77").Value = 3911.111
$ws.Range("I77").Value = 3837.5
$ws.Range("J77").Value = 4500
$ws.Range("K77").Value = 19187.5
$ws.Range("L77").Value = 22500
$ws.Range("M77").Value = -14507.5
$ws.Range("N77").Value = -31860

$ws.Range("H79").Value = 3068.7317
$ws.Range("I79").Value = 3097.6572
$ws.Range("J79").Value = 2900
$ws.Range("K79").Value = 3097.6572
$ws.Range("L79").Value = 2900
$ws.Range("M79").Value = -2005.6572
$ws.Range("N79").Value = -5084

$ws.Range("H125").Value = 1637.5
$ws.Range("I125").Value = 1650
$ws.Range("J125").Value = 1600
$ws.Range("K125").Value = 14850
$ws.Range("L125").Value = 14400
$ws.Range("M125").Value = -12390
$ws.Range("N125").Value = -19320

$ws.Range("H129").Value = 1657.3334
$ws.Range("J129").Value = 2869.25
$ws.Range("L129").Value = 8607.75
$ws.Range("N129").Value = -18607.75

$ws.Range("H138").Value = 2584.3333
$ws.Range("I138").Value = 2013.742
$ws.Range("J138").Value = 3194.276
$ws.Range("K138").Value = 6041.226
$ws.Range("L138").Value = 9582.828
$ws.Range("M138").Value = -901.2259999999997
$ws.Range("N138").Value = -19862.828

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 865.55554
$ws.Range("I45").Value = 838.8889
$ws.Range("J45").Value = 892.2222
$ws.Range("K45").Value = 838.8889
$ws.Range("L45").Value = 892.2222
$ws.Range("M45").Value = -461.8889
$ws.Range("N45").Value = -1646.2222

$ws.Range("H63").Value = 3777.6316
$ws.Range("I63").Value = 1510.7142
$ws.Range("J63").Value = 5100
$ws.Range("K63").Value = 1510.7142
$ws.Range("L63").Value = 5100
$ws.Range("M63").Value = -824.7141999999999
$ws.Range("N63").Value = -6472

$ws.Range("H66").Value = 3777.6316
$ws.Range("I66").Value = 1510.7142
$ws.Range("J66").Value = 5100
$ws.Range("K66").Value = 7553.571
$ws.Range("L66").Value = 25500
$ws.Range("M66").Value = -4121.571
$ws.Range("N66").Value = -32364

$ws.Range("H97").Value = 420.7143
$ws.Range("I97").Value = 375.30768
$ws.Range("J97").Value = 1011
$ws.Range("K97").Value = 375.30768
$ws.Range("L97").Value = 1011
$ws.Range("M97").Value = 120.69232
$ws.Range("N97").Value = -2003

$ws.Range("H110").Value = 1098.5
$ws.Range("I110").Value = 1008
$ws.Range("J110").Value = 1370
$ws.Range("K110").Value = 1008
$ws.Range("L110").Value = 1370
$ws.Range("M110").Value = 1037
$ws.Range("N110").Value = -5460

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2808.3333
$ws.Range("I62").Value = 2475
$ws.Range("K62").Value = 2475
$ws.Range("M62").Value = -1851

$ws.Range("H65").Value = 2808.3333
$ws.Range("I65").Value = 2475
$ws.Range("K65").Value = 12375
$ws.Range("M65").Value = -9255

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1284094.2
$ws.Range("J131").Value = 1588966.8
$ws.Range("L131").Value = 4766900.4
$ws.Range("N131").Value = -4776980.4

$ws.Range("H133").Value = 7291.087
$ws.Range("I133").Value = 4631.6665
$ws.Range("J133").Value = 8229.706
$ws.Range("K133").Value = 13894.9995
$ws.Range("L133").Value = 24689.118
$ws.Range("M133").Value = -8834.999500000002
$ws.Range("N133").Value = -34809.118

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 65.583336
$ws.Range("I2").Value = 37.875
$ws.Range("J2").Value = 121
$ws.Range("K2").Value = 37.875
$ws.Range("L2").Value = 121
$ws.Range("M2").Value = 75.125
$ws.Range("N2").Value = -347

$ws.Range("H11").Value = 2125

$ws.Range("H43").Value = 21571.428
$ws.Range("J43").Value = 30000
$ws.Range("L43").Value = 30000
$ws.Range("N43").Value = -30302

$ws.Range("H46").Value = 15616.667
$ws.Range("J46").Value = 15616.667
$ws.Range("L46").Value = 15616.667
$ws.Range("N46").Value = -15928.667

$ws.Range("H57").Value = 16294.75
$ws.Range("J57").Value = 16294.75
$ws.Range("L57").Value = 16294.75
$ws.Range("N57").Value = -17934.75

$ws.Range("H70").Value = 5609.8184
$ws.Range("I70").Value = 5545.3335
$ws.Range("J70").Value = 5900
$ws.Range("K70").Value = 5545.3335
$ws.Range("L70").Value = 5900
$ws.Range("M70").Value = -5275.3335
$ws.Range("N70").Value = -6440

$ws.Range("H73").Value = 5609.8184
$ws.Range("I73").Value = 5545.3335
$ws.Range("J73").Value = 5900
$ws.Range("K73").Value = 5545.3335
$ws.Range("L73").Value = 5900
$ws.Range("M73").Value = -4609.3335
$ws.Range("N73").Value = -7772

$ws.Range("H132").Value = 2817.074
$ws.Range("I132").Value = 2561.0557
$ws.Range("J132").Value = 3329.111
$ws.Range("K132").Value = 7683.1671
$ws.Range("L132").Value = 9987.332999999999
$ws.Range("M132").Value = -5153.1671
$ws.Range("N132").Value = -15047.333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 683.5
$ws.Range("I22").Value = 1000.5
$ws.Range("J22").Value = 525
$ws.Range("K22").Value = 1000.5
$ws.Range("L22").Value = 525
$ws.Range("M22").Value = -705.5
$ws.Range("N22").Value = -1115

$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()

$ws.Range("H27").Value = 683.5
$ws.Range("I27").Value = 1000.5
$ws.Range("J27").Value = 525
$ws.Range("K27").Value = 1000.5
$ws.Range("L27").Value = 525
$ws.Range("M27").Value = -893.5
$ws.Range("N27").Value = -739

$ws.Range("H94").Value = 30000
$ws.Range("J94").Value = 30000
$ws.Range("L94").Value = 30000
$ws.Range("N94").Value = -31352

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 3000
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").ClearContents()

$ws.Range("H54").Value = 15864.5
$ws.Range("J54").Value = 15864.5
$ws.Range("L54").Value = 15864.5
$ws.Range("N54").Value = -16904.5

$ws.Range("H81").Value = 3305.4666
$ws.Range("I81").Value = 1285.7142
$ws.Range("K81").Value = 2571.4284
$ws.Range("M81").Value = -1510.4284

$ws.Range("H84").Value = 3305.4666
$ws.Range("I84").Value = 1285.7142
$ws.Range("K84").Value = 12857.142
$ws.Range("M84").Value = -7553.142

$ws.Range("H136").Value = 9455.695
$ws.Range("I136").Value = 9749.137000000001
$ws.Range("K136").Value = 29247.411
$ws.Range("M136").Value = -26697.411
